$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.3526108804440327
$ws.Range("D2").Value = 0.144026428343949
$ws.Range("E2").Value = 0.1420455687202589
$ws.Range("F2").Value = 1.676978675952185
$ws.Range("G2").Value = 1.080943217832854
$ws.Range("H2").Value = 1.042518565243768
$ws.Range("J2").Value = 0.1734338221595539
$ws.Range("L2").Value = 0.13112348263809
$ws.Range("M2").Value = 1.313372010599664
$ws.Range("N2").Value = 1.834185285127148
$ws.Range("O2").Value = 4.312824721049111
$ws.Range("C3").Value = 0.3535761488437004
$ws.Range("D3").Value = 0.1450179076528251
$ws.Range("E3").Value = 0.143695570425499
$ws.Range("F3").Value = 1.666795739030661
$ws.Range("G3").Value = 1.064826376541674
$ws.Range("H3").Value = 1.040150824912715
$ws.Range("J3").Value = 0.1759001100428419
$ws.Range("L3").Value = 0.1319447406444088
$ws.Range("M3").Value = 1.204970989054431
$ws.Range("N3").Value = 1.689828024268877
$ws.Range("O3").Value = 4.27328729660843
$ws.Range("C4").Value = 0.3543545774805636
$ws.Range("D4").Value = 0.1456744167182471
$ws.Range("E4").Value = 0.1447691074980675
$ws.Range("F4").Value = 1.661507740557155
$ws.Range("G4").Value = 1.05563305011961
$ws.Range("H4").Value = 1.039218807273642
$ws.Range("J4").Value = 0.1774986309512512
$ws.Range("L4").Value = 0.1324822251331739
$ws.Range("M4").Value = 1.138260951558081
$ws.Range("N4").Value = 1.601396256336074
$ws.Range("O4").Value = 4.251577047033095
$ws.Range("C5").Value = 0.3547184730546888
$ws.Range("D5").Value = 0.1459539630109958
$ws.Range("E5").Value = 0.1452217696669855
$ws.Range("F5").Value = 1.65959512234943
$ws.Range("G5").Value = 1.052063025264303
$ws.Range("H5").Value = 1.038970192343456
$ws.Range("J5").Value = 0.178171190577757
$ws.Range("L5").Value = 0.1327096285888629
$ws.Range("M5").Value = 1.111040843189642
$ws.Range("N5").Value = 1.565415292377224
$ws.Range("O5").Value = 4.243374445446108
$ws.Range("C6").Value = 0.3547817161373246
$ws.Range("D6").Value = 0.1460011073335004
$ws.Range("E6").Value = 0.1452978508943872
$ws.Range("F6").Value = 1.659292160745096
$ws.Range("G6").Value = 1.05148086615128
$ws.Range("H6").Value = 1.038936833112885
$ws.Range("J6").Value = 0.1782841452032002
$ws.Range("L6").Value = 0.1327478950833676
$ws.Range("M6").Value = 1.106518921028581
$ws.Range("N6").Value = 1.559444173895713
$ws.Range("O6").Value = 4.242051315385083
$ws.Range("C7").Value = 0.3543592961368631
$ws.Range("D7").Value = 0.1456781381219088
$ws.Range("E7").Value = 0.1447751507776633
$ws.Range("F7").Value = 1.661480965607851
$ws.Range("G7").Value = 1.055584189979584
$ws.Range("H7").Value = 1.039214923210551
$ws.Range("J7").Value = 0.1775076157275155
$ws.Range("L7").Value = 0.1324852580451559
$ws.Range("M7").Value = 1.137893990521746
$ws.Range("N7").Value = 1.600910773136775
$ws.Range("O7").Value = 4.251463815275827
$ws.Range("C8").Value = 0.3529051345945078
$ws.Range("D8").Value = 0.1443583885089232
$ws.Range("E8").Value = 0.1426019423416083
$ws.Range("F8").Value = 1.673267308302755
$ws.Range("G8").Value = 1.075240134167117
$ws.Range("H8").Value = 1.041593857838507
$ws.Range("J8").Value = 0.1742667050750271
$ws.Range("L8").Value = 0.1313997669142815
$ws.Range("M8").Value = 1.276028420190798
$ws.Range("N8").Value = 1.784371078146592
$ws.Range("O8").Value = 4.298659212458261
$ws.Range("C9").Value = 0.3515289957243368
$ws.Range("D9").Value = 0.1421487001200674
$ws.Range("E9").Value = 0.1388200844130978
$ws.Range("F9").Value = 1.704045332490722
$ws.Range("G9").Value = 1.11937696905548
$ws.Range("H9").Value = 1.050401493878468
$ws.Range("J9").Value = 0.1685805188267242
$ws.Range("L9").Value = 0.129533943900114
$ws.Range("M9").Value = 1.545597036297991
$ws.Range("N9").Value = 2.145584680200557
$ws.Range("O9").Value = 4.411611586767208
$ws.Range("C10").Value = 0.3514201554665419
$ws.Range("D10").Value = 0.140755295079483
$ws.Range("E10").Value = 0.1363344045576493
$ws.Range("F10").Value = 1.731354270316118
$ws.Range("G10").Value = 1.155242780855161
$ws.Range("H10").Value = 1.05940362529617
$ws.Range("J10").Value = 0.1648121650986543
$ws.Range("L10").Value = 0.1283221969764021
$ws.Range("M10").Value = 1.742718355927494
$ws.Range("N10").Value = 2.411650076792796
$ws.Range("O10").Value = 4.507111198039581
$ws.Range("C11").Value = 0.3515671710755583
$ws.Range("D11").Value = 0.1401712293721182
$ws.Range("E11").Value = 0.1352672687811207
$ws.Range("F11").Value = 1.744802999042065
$ws.Range("G11").Value = 1.17231242861908
$ws.Range("H11").Value = 1.064049966642813
$ws.Range("J11").Value = 0.1631869937714026
$ws.Range("L11").Value = 0.1278052402267491
$ws.Range("M11").Value = 1.832165532483288
$ws.Range("N11").Value = 2.532798010411398
$ws.Range("O11").Value = 4.553291016980779
$ws.Range("C12").Value = 0.3516511458907559
$ws.Range("D12").Value = 0.139957211197661
$ws.Range("E12").Value = 0.1348723288916389
$ws.Range("F12").Value = 1.750043526625703
$ws.Range("G12").Value = 1.178885140369061
$ws.Range("H12").Value = 1.065888755017795
$ws.Range("J12").Value = 0.1625844201443272
$ws.Range("L12").Value = 0.1276143920379837
$ws.Range("M12").Value = 1.866001998021858
$ws.Range("N12").Value = 2.578686026819071
$ws.Range("O12").Value = 4.571172727298972
$ws.Range("C13").Value = 0.3516318009367154
$ws.Range("D13").Value = 0.1400029857310408
$ws.Range("E13").Value = 0.1349569785916893
$ws.Range("F13").Value = 1.748908307893259
$ws.Range("G13").Value = 1.177464744240069
$ws.Range("H13").Value = 1.06548921039672
$ws.Range("J13").Value = 0.1627136233777042
$ws.Range("L13").Value = 0.1276552764009971
$ws.Range("M13").Value = 1.858716324349274
$ws.Range("N13").Value = 2.568802760504354
$ws.Range("O13").Value = 4.567304024920531
$ws.Range("C14").Value = 0.3515735122499848
$ws.Range("D14").Value = 0.1401534786119285
$ws.Range("E14").Value = 0.1352345931535179
$ws.Range("F14").Value = 1.745231176921166
$ws.Range("G14").Value = 1.172850986820805
$ws.Range("H14").Value = 1.064199654868958
$ws.Range("J14").Value = 0.1631371620722017
$ws.Range("L14").Value = 0.1277894406783773
$ws.Range("M14").Value = 1.834950002141127
$ws.Range("N14").Value = 2.536573031141927
$ws.Range("O14").Value = 4.554754246209427
$ws.Range("C15").Value = 0.3515414959303911
$ws.Range("D15").Value = 0.1402465914240807
$ws.Range("E15").Value = 0.1354058335692931
$ws.Range("F15").Value = 1.74299808255445
$ws.Range("G15").Value = 1.17003910791027
$ws.Range("H15").Value = 1.063420096273262
$ws.Range("J15").Value = 0.1633982653970856
$ws.Range("L15").Value = 0.1278722593932073
$ws.Range("M15").Value = 1.820387766121613
$ws.Range("N15").Value = 2.516832814637837
$ws.Range("O15").Value = 4.547118535679772
$ws.Range("C16").Value = 0.3514145059394167
$ws.Range("D16").Value = 0.1407944667161871
$ws.Range("E16").Value = 0.1364054256200617
$ws.Range("F16").Value = 1.730496029283216
$ws.Range("G16").Value = 1.15414244515398
$ws.Range("H16").Value = 1.059111072634948
$ws.Range("J16").Value = 0.1649201692118378
$ws.Range("L16").Value = 0.1283566695789631
$ws.Range("M16").Value = 1.736868003509571
$ws.Range("N16").Value = 2.403734680770071
$ws.Range("O16").Value = 4.50414838180717
$ws.Range("C17").Value = 0.3513869674881676
$ws.Range("D17").Value = 0.1411433210920023
$ws.Range("E17").Value = 0.137034945493296
$ws.Range("F17").Value = 1.723089346337218
$ws.Range("G17").Value = 1.144583763528402
$ws.Range("H17").Value = 1.056608846548642
$ws.Range("J17").Value = 0.1658766462089591
$ws.Range("L17").Value = 0.1286626060317673
$ws.Range("M17").Value = 1.685571756057243
$ws.Range("N17").Value = 2.334378555150352
$ws.Range("O17").Value = 4.478489123378438
$ws.Range("C18").Value = 0.3513896236381839
$ws.Range("D18").Value = 0.1413486608168384
$ws.Range("E18").Value = 0.1374030153451375
$ws.Range("F18").Value = 1.718925750128733
$ws.Range("G18").Value = 1.139156821092882
$ws.Range("H18").Value = 1.055221513001953
$ws.Range("J18").Value = 0.1664351682172587
$ws.Range("L18").Value = 0.1288417994893294
$ws.Range("M18").Value = 1.656046576288148
$ws.Range("N18").Value = 2.294497683618829
$ws.Range("O18").Value = 4.463988120995225
$ws.Range("C19").Value = 0.3513936982925827
$ws.Range("D19").Value = 0.1414189906730208
$ws.Range("E19").Value = 0.1375286652677783
$ws.Range("F19").Value = 1.717532598846219
$ws.Range("G19").Value = 1.137331527665623
$ws.Range("H19").Value = 1.054760694648593
$ws.Range("J19").Value = 0.1666257126785773
$ws.Range("L19").Value = 0.1289030260377757
$ws.Range("M19").Value = 1.646046353384165
$ws.Range("N19").Value = 2.280996718963991
$ws.Range("O19").Value = 4.459122533500249
$ws.Range("C20").Value = 0.3513879844749681
$ws.Range("D20").Value = 0.1411056998143749
$ws.Range("E20").Value = 0.1369673123296709
$ws.Range("F20").Value = 1.72386780763614
$ws.Range("G20").Value = 1.14559395492131
$ws.Range("H20").Value = 1.056869842935669
$ws.Range("J20").Value = 0.1657739598671579
$ws.Range("L20").Value = 0.1286297047200939
$ws.Range("M20").Value = 1.691034517533609
$ws.Range("N20").Value = 2.34176053310108
$ws.Range("O20").Value = 4.4811939359127
$ws.Range("C21").Value = 0.3515898645686377
$ws.Range("D21").Value = 0.1401090810764885
$ws.Range("E21").Value = 0.135152802282037
$ws.Range("F21").Value = 1.746307226287513
$ws.Range("G21").Value = 1.174203203307201
$ws.Range("H21").Value = 1.064576275433467
$ws.Range("J21").Value = 0.1630124096416861
$ws.Range("L21").Value = 0.1277499001938693
$ws.Range("M21").Value = 1.841931727150097
$ws.Range("N21").Value = 2.546039395052844
$ws.Range("O21").Value = 4.558429708912627
$ws.Range("C22").Value = 0.3518867975492554
$ws.Range("D22").Value = 0.1394994339779529
$ws.Range("E22").Value = 0.1340203204498431
$ws.Range("F22").Value = 1.761834225716427
$ws.Range("G22").Value = 1.193535380062798
$ws.Range("H22").Value = 1.070075236873947
$ws.Range("J22").Value = 0.1612824549919587
$ws.Range("L22").Value = 0.1272035211619098
$ws.Range("M22").Value = 1.94034500484716
$ws.Range("N22").Value = 2.679614937623569
$ws.Range("O22").Value = 4.611207126307647
$ws.Range("C23").Value = 0.3517132078045933
$ws.Range("D23").Value = 0.1398210006158251
$ws.Range("E23").Value = 0.1346198561190173
$ws.Range("F23").Value = 1.753468248603738
$ws.Range("G23").Value = 1.183159266364612
$ws.Range("H23").Value = 1.067098013252263
$ws.Range("J23").Value = 0.1621989013602825
$ws.Range("L23").Value = 0.1274925201791319
$ws.Range("M23").Value = 1.8878398875286
$ws.Range("N23").Value = 2.608318417713463
$ws.Range("O23").Value = 4.582828133570729
$ws.Range("C24").Value = 0.3513874671098165
$ws.Range("D24").Value = 0.1411226935084144
$ws.Range("E24").Value = 0.1369978701253158
$ws.Range("F24").Value = 1.723515570608555
$ws.Range("G24").Value = 1.145137034133683
$ws.Range("H24").Value = 1.056751686887225
$ws.Range("J24").Value = 0.1658203574799675
$ws.Range("L24").Value = 0.1286445691015068
$ws.Range("M24").Value = 1.688564909626649
$ws.Range("N24").Value = 2.338423162681011
$ws.Range("O24").Value = 4.479970308835107
$ws.Range("C25").Value = 0.3517430381108113
$ws.Range("D25").Value = 0.1427060371992397
$ws.Range("E25").Value = 0.1397917753362279
$ws.Range("F25").Value = 1.694896043801592
$ws.Range("G25").Value = 1.106835103166674
$ws.Range("H25").Value = 1.047574717685251
$ws.Range("J25").Value = 0.1700470029199002
$ws.Range("L25").Value = 0.1300106810094785
$ws.Range("M25").Value = 1.472825668018615
$ws.Range("N25").Value = 2.047732773542464
$ws.Range("O25").Value = 4.378862964720952
